# Generate Report for Handback
# This script updates the localization-status workbook to reflect a completed
# handback: the overview status text is updated, and the zh-cn / de-de detail
# sheets get their "Latest Target File", "Latest Handback File" and
# "Latest Handback DateTime" columns filled in (with a hyperlink on the target
# file name), for both tracked source files.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdFile1 = "503536e2-20fd-45df-a15a-f4fcab9f39e3.md"
$mdFile2 = "92a6c502-9cba-4da6-b279-a5e3cbc954a7.md"
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93f04f9f90c6b187952aa96e75ae580e3f23a157/e2e/503536e2-20fd-45df-a15a-f4fcab9f39e3.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/93f04f9f90c6b187952aa96e75ae580e3f23a157/e2e/92a6c502-9cba-4da6-b279-a5e3cbc954a7.md"

# ---------------------------------------------------------------------------
# Overview sheet: mark both tracked files as handed back (zh-cn & de-de cols)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# widen the zh-cn / de-de status columns so the longer text fits
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("I2").Value = $mdFile1
$wsZhCn.Range("J2").Value = "503536e2-20fd-45df-a15a-f4fcab9f39e3.c80a21061b4aedf9a65b75648d39be0efcba674a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-20 14:34:48"

$wsZhCn.Range("I3").Value = $mdFile2
$wsZhCn.Range("J3").Value = "92a6c502-9cba-4da6-b279-a5e3cbc954a7.8cf542d7008303cc24bbc864491bab1555247e42.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-20 14:34:48"

$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl1, "", "", $mdFile1)
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 15570276

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("I2").Value = $mdFile1
$wsDeDe.Range("J2").Value = "503536e2-20fd-45df-a15a-f4fcab9f39e3.c80a21061b4aedf9a65b75648d39be0efcba674a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-20 14:34:53"

$wsDeDe.Range("I3").Value = $mdFile2
$wsDeDe.Range("J3").Value = "92a6c502-9cba-4da6-b279-a5e3cbc954a7.8cf542d7008303cc24bbc864491bab1555247e42.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-20 14:34:53"

$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl1, "", "", $mdFile1)
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 15570276

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
